# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K") previously held a raw strike-count
# value ("Strike#"). The data pipeline was re-run against the underlying
# option chain so that column G now reports the (reduced) strike offset
# "K" for each trade row. Only the values in column G change; everything
# else on the sheet (dates, TB/PC/dS0/dSF/IP/I0/IF, headers, styles)
# stays exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values for rows 2-43, in row order.
$kValues = @(
    0,0,1,1,2,0,1,2,1,1,
    0,2,1,0,3,0,1,2,0,0,
    2,2,1,1,2,1,1,0,2,0,
    2,0,2,1,1,0,0,1,1,1,
    1,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
